# "Add Group set/reset buttons in dataloader page"
#
# Timeline tracker update: log a new line item for the work done adding the
# Group set/reset buttons on the dataloader page. Row 21 was the next blank
# entry row in the table, so it gets the task description, hours spent, and
# the completion date; the existing shared formula in column D (hours *
# hourly rate) then recalculates automatically, which also rolls the new
# total up into the ИТОГО (SUM) cell at the bottom of the sheet. Finally,
# move the sheet's active selection down to C22, the next empty row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A21").Value = "Работа по созданию функционала загрузки данных (Первый альфа вариант)"
$ws.Range("B21").Value = 3

# Pull C21's date formatting from C20 (same column style already used for
# every other "Дата выполнения" cell) before writing the serial date value,
# so it renders/saves with the existing date number format instead of
# picking up a brand-new one.
$ws.Range("C20").Copy()
$ws.Range("C21").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C21").Value = 43555   # 2019-03-31

# D21 keeps its pre-existing shared formula (=B21*$B$1 -> 3 * 800 = 2400);
# no explicit write needed, recalculation after the script refreshes it
# along with the D38 SUM(D4:D37) grand total (28800 -> 31200).

$ws.Range("C22").Select()
